$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add 5 new sheets after "Day 12" (the current last sheet).
#
# Creation order matters: the target sheetId / r:id allocation in the
# workbook shows the sheets were created in the order
#   Day 13 emp -> Day 14 -> Day 13 dept -> Day 15 -> Day 16
# with "Day 13 dept" then moved so it sits right after "Day 13 emp".
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsEmp = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsEmp.Name = "Day 13 emp"

$wsDay14 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsEmp)
$wsDay14.Name = "Day 14"

$wsDept = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsDay14)
$wsDept.Name = "Day 13 dept"

# Move "Day 13 dept" so it sits right after "Day 13 emp" (before "Day 14")
$wsDept.Move([System.Reflection.Missing]::Value, $wsEmp)

# NOTE: after the Move above, sheet positions shift - reference sheets by
# re-fetching them by name so later steps are unambiguous.
$wsEmp   = $wb.Worksheets.Item("Day 13 emp")
$wsDept  = $wb.Worksheets.Item("Day 13 dept")
$wsDay14 = $wb.Worksheets.Item("Day 14")

$wsDay15 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsDay14)
$wsDay15.Name = "Day 15"

$wsDay16 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsDay15)
$wsDay16.Name = "Day 16"

# ------------------------------------------------------------------
# Day 13 emp: employees table (id, name, salary, departmentId)
# ------------------------------------------------------------------
$wsEmp.Cells.Item(1,1).Value = "id"
$wsEmp.Cells.Item(1,2).Value = "name"
$wsEmp.Cells.Item(1,3).Value = "salary"
$wsEmp.Cells.Item(1,4).Value = "departmentId"

$wsEmp.Cells.Item(2,1).Value = 1
$wsEmp.Cells.Item(2,2).Value = "Joe"
$wsEmp.Cells.Item(2,3).Value = 70000
$wsEmp.Cells.Item(2,4).Value = 1

$wsEmp.Cells.Item(3,1).Value = 2
$wsEmp.Cells.Item(3,2).Value = "Jim"
$wsEmp.Cells.Item(3,3).Value = 90000
$wsEmp.Cells.Item(3,4).Value = 1

$wsEmp.Cells.Item(4,1).Value = 3
$wsEmp.Cells.Item(4,2).Value = "Henry"
$wsEmp.Cells.Item(4,3).Value = 80000
$wsEmp.Cells.Item(4,4).Value = 2

$wsEmp.Cells.Item(5,1).Value = 4
$wsEmp.Cells.Item(5,2).Value = "Sam"
$wsEmp.Cells.Item(5,3).Value = 60000
$wsEmp.Cells.Item(5,4).Value = 2

$wsEmp.Cells.Item(6,1).Value = 5
$wsEmp.Cells.Item(6,2).Value = "Max"
$wsEmp.Cells.Item(6,3).Value = 90000
$wsEmp.Cells.Item(6,4).Value = 1

$wsEmp.Columns.Item(4).AutoFit()
$wsEmp.Range("D7").Select()

# ------------------------------------------------------------------
# Day 13 dept: departments table (id, name)
# ------------------------------------------------------------------
$wsDept.Cells.Item(1,1).Value = "id"
$wsDept.Cells.Item(1,2).Value = "name"

$wsDept.Cells.Item(2,1).Value = 1
$wsDept.Cells.Item(2,2).Value = "IT"

$wsDept.Cells.Item(3,1).Value = 2
$wsDept.Cells.Item(3,2).Value = "Sales"

$wsDept.Range("B4").Select()

# ------------------------------------------------------------------
# Day 14: scores table (id, score)
# ------------------------------------------------------------------
$wsDay14.Cells.Item(1,1).Value = "id"
$wsDay14.Cells.Item(1,2).Value = "score"

$wsDay14.Cells.Item(2,1).Value = 1
$wsDay14.Cells.Item(2,2).Value = 3.5

$wsDay14.Cells.Item(3,1).Value = 2
$wsDay14.Cells.Item(3,2).Value = 3.65

$wsDay14.Cells.Item(4,1).Value = 3
$wsDay14.Cells.Item(4,2).Value = 4

$wsDay14.Cells.Item(5,1).Value = 4
$wsDay14.Cells.Item(5,2).Value = 3.85

$wsDay14.Cells.Item(6,1).Value = 5
$wsDay14.Cells.Item(6,2).Value = 4

$wsDay14.Cells.Item(7,1).Value = 6
$wsDay14.Cells.Item(7,2).Value = 3.65

$wsDay14.Range("H18").Select()

# ------------------------------------------------------------------
# Day 15: user emails table (id, email)
# ------------------------------------------------------------------
$wsDay15.Cells.Item(1,1).Value = "id"
$wsDay15.Cells.Item(1,2).Value = "email"

$wsDay15.Cells.Item(2,1).Value = 1
$wsDay15.Cells.Item(2,2).Value = "john@example.com"

$wsDay15.Cells.Item(3,1).Value = 2
$wsDay15.Cells.Item(3,2).Value = "bob@example.com"

$wsDay15.Cells.Item(4,1).Value = 3
$wsDay15.Cells.Item(4,2).Value = "john@example.com"

$wsDay15.Columns.Item(2).AutoFit()
$wsDay15.Range("G13").Select()

# ------------------------------------------------------------------
# Day 16: store prices table (product_id, store1, store2, store3)
# ------------------------------------------------------------------
$wsDay16.Cells.Item(1,1).Value = "product_id"
$wsDay16.Cells.Item(1,2).Value = "store1"
$wsDay16.Cells.Item(1,3).Value = "store2"
$wsDay16.Cells.Item(1,4).Value = "store3"

$wsDay16.Cells.Item(2,1).Value = 0
$wsDay16.Cells.Item(2,2).Value = 95
$wsDay16.Cells.Item(2,3).Value = 100
$wsDay16.Cells.Item(2,4).Value = 105

$wsDay16.Cells.Item(3,1).Value = 1
$wsDay16.Cells.Item(3,2).Value = 70
$wsDay16.Cells.Item(3,4).Value = 80

$wsDay16.Columns.Item(1).AutoFit()
$wsDay16.Range("L13").Select()

# Day 16 is the sheet that ends up active/selected in the target workbook.
$wsDay16.Select()
